$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates ---
$ws.Range("A8").Value = "Volume 32   Number  24"
$ws.Range("C9").Value = "Report Covering the Week  6/9/2025  Through  6/15/2025"

# --- Type-changing cells: number -> text (use NumberFormat to force text, then restore original style) ---
$ws.Range("F15").NumberFormat = "@"
$ws.Range("F15").Value = "0"
$ws.Range("C15").Copy()
$ws.Range("F15").PasteSpecial(-4122)
$ws.Range("C20").NumberFormat = "@"
$ws.Range("C20").Value = "0"
$ws.Range("C15").Copy()
$ws.Range("C20").PasteSpecial(-4122)
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0"
$ws.Range("C15").Copy()
$ws.Range("D20").PasteSpecial(-4122)
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "***.*"
$ws.Range("C15").Copy()
$ws.Range("E20").PasteSpecial(-4122)
$ws.Range("C22").NumberFormat = "@"
$ws.Range("C22").Value = "0"
$ws.Range("C15").Copy()
$ws.Range("C22").PasteSpecial(-4122)
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0"
$ws.Range("C15").Copy()
$ws.Range("D23").PasteSpecial(-4122)
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "***.*"
$ws.Range("C15").Copy()
$ws.Range("E23").PasteSpecial(-4122)
$ws.Range("F27").NumberFormat = "@"
$ws.Range("F27").Value = "0"
$ws.Range("C15").Copy()
$ws.Range("F27").PasteSpecial(-4122)

# --- Type-changing cells: text -> number (copy numeric format first, then set value) ---
$ws.Range("C16").Copy()
$ws.Range("D22").PasteSpecial(-4122)
$ws.Range("D22").Value = 2
$ws.Range("E16").Copy()
$ws.Range("E22").PasteSpecial(-4122)
$ws.Range("E22").Value = -100
$ws.Range("C16").Copy()
$ws.Range("D27").PasteSpecial(-4122)
$ws.Range("D27").Value = 2
$ws.Range("E16").Copy()
$ws.Range("E27").PasteSpecial(-4122)
$ws.Range("E27").Value = -100
$ws.Range("C16").Copy()
$ws.Range("C28").PasteSpecial(-4122)
$ws.Range("C28").Value = 5
$ws.Range("C16").Copy()
$ws.Range("D28").PasteSpecial(-4122)
$ws.Range("D28").Value = 1
$ws.Range("E16").Copy()
$ws.Range("E28").PasteSpecial(-4122)
$ws.Range("E28").Value = 400

# --- Plain numeric value updates ---
$ws.Range("G15").Value = 1
$ws.Range("H15").Value = -100
$ws.Range("L15").Value = 83.333333333333
$ws.Range("N15").Value = 37.5
$ws.Range("C16").Value = 4
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 300
$ws.Range("F16").Value = 18
$ws.Range("G16").Value = 13
$ws.Range("H16").Value = 38.461538461538
$ws.Range("I16").Value = 85
$ws.Range("J16").Value = 65
$ws.Range("K16").Value = 30.769230769230
$ws.Range("L16").Value = 13.333333333333
$ws.Range("M16").Value = 8.974358974358
$ws.Range("N16").Value = -82.067510548523
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = 40
$ws.Range("F17").Value = 22
$ws.Range("G17").Value = 14
$ws.Range("H17").Value = 57.142857142857
$ws.Range("I17").Value = 125
$ws.Range("J17").Value = 94
$ws.Range("K17").Value = 32.978723404255
$ws.Range("L17").Value = 31.578947368421
$ws.Range("M17").Value = 115.51724137931
$ws.Range("N17").Value = 6.837606837606
$ws.Range("C18").Value = 5
$ws.Range("D18").Value = 4
$ws.Range("E18").Value = 25
$ws.Range("F18").Value = 27
$ws.Range("G18").Value = 16
$ws.Range("H18").Value = 68.75
$ws.Range("I18").Value = 166
$ws.Range("J18").Value = 88
$ws.Range("K18").Value = 88.636363636363
$ws.Range("L18").Value = 48.214285714285
$ws.Range("M18").Value = 22.962962962963
$ws.Range("N18").Value = -77.597840755735
$ws.Range("D19").Value = 16
$ws.Range("E19").Value = 12.5
$ws.Range("F19").Value = 86
$ws.Range("H19").Value = 24.637681159420
$ws.Range("I19").Value = 491
$ws.Range("J19").Value = 436
$ws.Range("K19").Value = 12.614678899082
$ws.Range("L19").Value = 6.739130434782
$ws.Range("M19").Value = -23.757763975155
$ws.Range("N19").Value = -58.319185059422
$ws.Range("F20").Value = 8
$ws.Range("G20").Value = 3
$ws.Range("H20").Value = 166.666666666667
$ws.Range("N20").Value = -96.092362344582
$ws.Range("C21").Value = 34
$ws.Range("D21").Value = 26
$ws.Range("E21").Value = 30.769230769230
$ws.Range("F21").Value = 161
$ws.Range("G21").Value = 116
$ws.Range("H21").Value = 38.793103448275
$ws.Range("I21").Value = 900
$ws.Range("J21").Value = 703
$ws.Range("K21").Value = 28.022759601707
$ws.Range("L21").Value = 15.384615384615
$ws.Range("M21").Value = -3.743315508021
$ws.Range("N21").Value = -70.826580226904
$ws.Range("F22").Value = 3
$ws.Range("G22").Value = 3
$ws.Range("H22").Value = 0
$ws.Range("J22").Value = 28
$ws.Range("K22").Value = 17.857142857142
$ws.Range("L22").Value = -13.157894736842
$ws.Range("M22").Value = -5.714285714285
$ws.Range("C24").Value = 74
$ws.Range("D24").Value = 42
$ws.Range("E24").Value = 76.190476190476
$ws.Range("F24").Value = 259
$ws.Range("G24").Value = 214
$ws.Range("H24").Value = 21.028037383177
$ws.Range("I24").Value = 1263
$ws.Range("J24").Value = 1352
$ws.Range("K24").Value = -6.582840236686
$ws.Range("L24").Value = 34.504792332268
$ws.Range("M24").Value = 57.677902621722
$ws.Range("C25").Value = 56
$ws.Range("D25").Value = 28
$ws.Range("E25").Value = 100
$ws.Range("G25").Value = 184
$ws.Range("H25").Value = 10.869565217391
$ws.Range("I25").Value = 1015
$ws.Range("J25").Value = 1151
$ws.Range("K25").Value = -11.815812337098
$ws.Range("L25").Value = 46.253602305475
$ws.Range("C26").Value = 12
$ws.Range("D26").Value = 16
$ws.Range("E26").Value = -25
$ws.Range("F26").Value = 39
$ws.Range("G26").Value = 59
$ws.Range("H26").Value = -33.898305084745
$ws.Range("I26").Value = 258
$ws.Range("J26").Value = 268
$ws.Range("K26").Value = -3.731343283582
$ws.Range("L26").Value = 18.894009216589
$ws.Range("M26").Value = 35.789473684210
$ws.Range("G27").Value = 3
$ws.Range("H27").Value = -100
$ws.Range("J27").Value = 8
$ws.Range("K27").Value = 87.5
$ws.Range("L27").Value = 25
$ws.Range("F28").Value = 9
$ws.Range("G28").Value = 9
$ws.Range("H28").Value = 0
$ws.Range("I28").Value = 57
$ws.Range("J28").Value = 51
$ws.Range("K28").Value = 11.764705882352
$ws.Range("L28").Value = 1.785714285714

$ws.Range("A1").Select()
